$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.781.93'
$ws.Range("E2").Value = '  +3.15%  '
$ws.Range("D3").Value = '4.027.62'
$ws.Range("E3").Value = '  +2.34%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '523.72'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.01'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.35%  '
$ws.Range("E7").Value = '  +0.78%  '
$ws.Range("E9").Value = '  +1.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.178'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000341'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.43'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +8.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.79'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.38%  '
$ws.Range("D14").Value = '4.670.58'
$ws.Range("E14").Value = '  +2.26%  '
$ws.Range("D15").Value = '4.060.36'
$ws.Range("E15").Value = '  +2.98%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.46'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +7.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.25'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.49%  '
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("E19").Value = '  -1.82%  '
$ws.Range("D20").Value = '71.705.56'
$ws.Range("E20").Value = '  +3.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '442.02'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.61'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +6.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '94.47'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +6.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.36'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.27'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.40%  '
$ws.Range("E26").Value = '  -1.02%  '
$ws.Range("E27").Value = '  +2.86%  '
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.62'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '699.06'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.132'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.93'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.96'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +13.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '67.41'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.67%  '
$ws.Range("D35").Value = '0.0₃0911'
$ws.Range("E35").Value = '  +5.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.445'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.86'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.157'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +5.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.57'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +18.80%  '
$ws.Range("E40").Value = '  +0.17%  '
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("E42").Value = '  +1.44%  '
$ws.Range("E43").Value = '  +0.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.82'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.52'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.146'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.21'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000279'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +18.05%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.20'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +6.09%  '
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0341'
$ws.Range("E51").Value = '  -5.04%  '
